$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at row 20, pushing rows 20-150 down to 21-151.
$ws.Rows.Item(20).Insert()

# The newly-inserted row copies formatting from the row above by default;
# re-apply the formatting (difficulty highlight colour) from the row that
# got pushed down to row 21 (the original row 20) so the new row keeps the
# same "medium" style that row used to have.
$ws.Range("A21:C21").Copy() | Out-Null
$ws.Range("A20:C20").PasteSpecial(-4122) | Out-Null

# Fill in the new LeetCode entry.
$ws.Range("A20").Value = "75. Sort Colors"
$ws.Range("B20").Value = "h"
$ws.Range("C20").Value = "bucket sort (use a bucket of size n (#of distinct values); 1 pass - 2 ptrs + partition"

# Match the saved view state (zoom level and active selection).
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 149
$ws.Range("C21").Select() | Out-Null
